# Refresh of the cryptos list (Price / Volume(1h) columns) plus a re-ranking
# of three coins (Injective Protocol, Rocket Pool ETH, Frax Share) whose rows
# changed order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 46-48 swapped which coin they describe (Coin name + Link).
$coinUpdates = @(
    @{ Row = 46; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj" },
    @{ Row = 47; B = "RocketPoolETH";     C = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth" },
    @{ Row = 48; B = "FraxShare";         C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" }
)

foreach ($u in $coinUpdates) {
    $ws.Range("B$($u.Row)").Value = $u.B
    $ws.Range("C$($u.Row)").Value = $u.C
}

# Updated Price (column D) and Volume(1h) (column E) values.
$priceUpdates = @(
    @{ Row = 2;  D = "34.336.62";  E = "  -0.59%  " },
    @{ Row = 3;  D = "1.786.75";   E = "  -2.82%  " },
    @{ Row = 4;  D = "0.998";      E = "  -0.09%  " },
    @{ Row = 5;  D = "224.50";     E = "  -3.13%  " },
    @{ Row = 6;  D = "0.552";      E = "  -4.30%  " },
    @{ Row = 7;  D = "0.998";      E = "  -0.11%  " },
    @{ Row = 8;  D = "33.31";      E = "  +4.96%  " },
    @{ Row = 9;  D = "0.281";      E = "  -2.74%  " },
    @{ Row = 10; D = "0.0658";     E = "  -3.65%  " },
    @{ Row = 11; D = "0.0932";     E = "  -0.40%  " },
    @{ Row = 12; D = "2.042.92";   E = "  -2.89%  " },
    @{ Row = 13; D = "11.11";      E = "  +7.04%  " },
    @{ Row = 14; D = "1.811.68";   E = "  -1.54%  " },
    @{ Row = 15; D = "0.631";      E = "  -3.34%  " },
    @{ Row = 16; D = "34.291.34";  E = "  -0.56%  " },
    @{ Row = 17; D = "4.23";       E = "  -2.02%  " },
    @{ Row = 18; D = "68.87";      E = "  -2.27%  " },
    @{ Row = 19; D = "254.97";     E = "  -2.33%  " },
    @{ Row = 20; D = "0.0₃0740";   E = "  -2.44%  " },
    @{ Row = 21; D = "0.998";      E = "  -0.18%  " },
    @{ Row = 22; D = "10.41";      E = "  -2.13%  " },
    @{ Row = 23; D = $null;        E = "  -5.42%  " },
    @{ Row = 24; D = "2.12";       E = "  -4.84%  " },
    @{ Row = 25; D = "157.27";     E = "  -0.82%  " },
    @{ Row = 26; D = "16.41";      E = "  -2.43%  " },
    @{ Row = 27; D = "7.01";       E = "  -2.79%  " },
    @{ Row = 28; D = $null;        E = "  -4.53%  " },
    @{ Row = 29; D = "0.999";      E = "  +0.04%  " },
    @{ Row = 30; D = "3.78";       E = "  -2.71%  " },
    @{ Row = 31; D = "0.0514";     E = "  -3.11%  " },
    @{ Row = 32; D = $null;        E = "  -2.84%  " },
    @{ Row = 33; D = $null;        E = "  -0.97%  " },
    @{ Row = 34; D = "1.88";       E = "  +3.63%  " },
    @{ Row = 35; D = "1.447.34";   E = "  -7.31%  " },
    @{ Row = 36; D = $null;        E = "  -2.51%  " },
    @{ Row = 37; D = "0.0188";     E = "  -1.88%  " },
    @{ Row = 38; D = "0.625";      E = "  -2.66%  " },
    @{ Row = 39; D = "2.86";       E = "  +1.08%  " },
    @{ Row = 40; D = "83.23";      E = "  -3.45%  " },
    @{ Row = 41; D = "2.34";       E = "  -0.28%  " },
    @{ Row = 42; D = "0.891";      E = "  -3.73%  " },
    @{ Row = 43; D = "2.07";       E = "  -3.90%  " },
    @{ Row = 44; D = "0.0506";     E = "  -3.98%  " },
    @{ Row = 45; D = $null;        E = "  -3.10%  " },
    @{ Row = 46; D = "12.35";      E = "  -1.41%  " },
    @{ Row = 47; D = "1.942.67";   E = "  -2.69%  " },
    @{ Row = 48; D = "5.81";       E = "  -2.01%  " },
    @{ Row = 49; D = "0.999";      E = "  -0.15%  " },
    @{ Row = 50; D = "99.01";      E = "  -0.07%  " },
    @{ Row = 51; D = "49.74";      E = "  -4.65%  " }
)

foreach ($u in $priceUpdates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$($u.Row)")
        # The Price column stores plain display strings (e.g. "34.336.62",
        # "0.998"). Values that parse as a plain number would otherwise be
        # silently turned into a numeric cell by Excel, so force text
        # storage first for anything that looks like a number.
        if ($u.D -match '^[+-]?[0-9]*\.?[0-9]+$') {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
